$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update status column (F) for rows 16, 21, 22 from "Por iniciar" to "Hecho"
$ws.Range("F16").Value = "Hecho"
$ws.Range("F21").Value = "Hecho"
$ws.Range("F22").Value = "Hecho"

# Record variable-expense values in column Z for rows 16, 21, 22
$ws.Range("Z16").Value = 1.5
$ws.Range("Z21").Value = 1
$ws.Range("Z22").Value = 0.5

# Update the active selection on the sheet to C5 (bottom-right frozen pane)
$ws.Activate()
$ws.Range("C5").Select()
